$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 Feature GA Results")
$ws.Range("A33:A131").FormulaR1C1 = "=R[-1]C+1"
$ws.Cells.Item(33, 2).Value = 62.9
$ws.Cells.Item(34, 2).Value = 61.4
$ws.Cells.Item(35, 2).Value = 70
$ws.Cells.Item(36, 2).Value = 64.3
$ws.Cells.Item(37, 2).Value = 61.4
$ws.Cells.Item(38, 2).Value = 70
$ws.Cells.Item(39, 2).Value = 58.6
$ws.Cells.Item(40, 2).Value = 67.099999999999994
$ws.Cells.Item(41, 2).Value = 62.9
$ws.Cells.Item(42, 2).Value = 68.599999999999994
$ws.Cells.Item(43, 2).Value = 61.4
$ws.Cells.Item(44, 2).Value = 64.3
$ws.Cells.Item(45, 2).Value = 62.9
$ws.Cells.Item(46, 2).Value = 64.3
$ws.Cells.Item(47, 2).Value = 64.3
$ws.Cells.Item(48, 2).Value = 64.3
$ws.Cells.Item(49, 2).Value = 64.3
$ws.Cells.Item(50, 2).Value = 70
$ws.Cells.Item(51, 2).Value = 62.9
$ws.Cells.Item(52, 2).Value = 77.099999999999994
$ws.Cells.Item(53, 2).Value = 62.9
$ws.Cells.Item(54, 2).Value = 67.099999999999994
$ws.Cells.Item(55, 2).Value = 58.6
$ws.Cells.Item(56, 2).Value = 62.9
$ws.Cells.Item(57, 2).Value = 65.7
$ws.Cells.Item(58, 2).Value = 61.4
$ws.Cells.Item(59, 2).Value = 64.3
$ws.Cells.Item(60, 2).Value = 60
$ws.Cells.Item(61, 2).Value = 64.3
$ws.Cells.Item(62, 2).Value = 62.9
$ws.Cells.Item(63, 2).Value = 62.9
$ws.Cells.Item(64, 2).Value = 62.9
$ws.Cells.Item(65, 2).Value = 60
$ws.Cells.Item(66, 2).Value = 60
$ws.Cells.Item(67, 2).Value = 57.1
$ws.Cells.Item(68, 2).Value = 62.9
$ws.Cells.Item(69, 2).Value = 64.3
$ws.Cells.Item(70, 2).Value = 61.4
$ws.Cells.Item(71, 2).Value = 55.7
$ws.Cells.Item(72, 2).Value = 61.4
$ws.Cells.Item(73, 2).Value = 61.4
$ws.Cells.Item(74, 2).Value = 64.3
$ws.Cells.Item(75, 2).Value = 62.9
$ws.Cells.Item(76, 2).Value = 65.7
$ws.Cells.Item(77, 2).Value = 62.9
$ws.Cells.Item(78, 2).Value = 60
$ws.Cells.Item(79, 2).Value = 60
$ws.Cells.Item(80, 2).Value = 65.7
$ws.Cells.Item(81, 2).Value = 62.9
$ws.Cells.Item(82, 2).Value = 68.599999999999994
$ws.Cells.Item(83, 2).Value = 57.1
$ws.Cells.Item(84, 2).Value = 61.4
$ws.Cells.Item(85, 2).Value = 62.9
$ws.Cells.Item(86, 2).Value = 68.599999999999994
$ws.Cells.Item(87, 2).Value = 60
$ws.Cells.Item(88, 2).Value = 57.1
$ws.Cells.Item(89, 2).Value = 58.6
$ws.Cells.Item(90, 2).Value = 68.599999999999994
$ws.Cells.Item(91, 2).Value = 67.099999999999994
$ws.Cells.Item(92, 2).Value = 67.099999999999994
$ws.Cells.Item(93, 2).Value = 65.7
$ws.Cells.Item(94, 2).Value = 60
$ws.Cells.Item(95, 2).Value = 61.4
$ws.Cells.Item(96, 2).Value = 62.9
$ws.Cells.Item(97, 2).Value = 58.6
$ws.Cells.Item(98, 2).Value = 61.4
$ws.Cells.Item(99, 2).Value = 54.3
$ws.Cells.Item(100, 2).Value = 65.7
$ws.Cells.Item(101, 2).Value = 57.1
$ws.Cells.Item(102, 2).Value = 64.3
$ws.Cells.Item(103, 2).Value = 58.6
$ws.Cells.Item(104, 2).Value = 67.099999999999994
$ws.Cells.Item(105, 2).Value = 57.1
$ws.Cells.Item(106, 2).Value = 65.7
$ws.Cells.Item(107, 2).Value = 61.4
$ws.Cells.Item(108, 2).Value = 57.1
$ws.Cells.Item(109, 2).Value = 61.4
$ws.Cells.Item(110, 2).Value = 64.3
$ws.Cells.Item(111, 2).Value = 64.3
$ws.Cells.Item(112, 2).Value = 67.099999999999994
$ws.Cells.Item(113, 2).Value = 58.6
$ws.Cells.Item(114, 2).Value = 67.099999999999994
$ws.Cells.Item(115, 2).Value = 58.6
$ws.Cells.Item(116, 2).Value = 62.9
$ws.Cells.Item(117, 2).Value = 65.7
$ws.Cells.Item(118, 2).Value = 62.9
$ws.Cells.Item(119, 2).Value = 62.9
$ws.Cells.Item(120, 2).Value = 64.3
$ws.Cells.Item(121, 2).Value = 57.1
$ws.Cells.Item(122, 2).Value = 67.099999999999994
$ws.Cells.Item(123, 2).Value = 67.099999999999994
$ws.Cells.Item(124, 2).Value = 61.4
$ws.Cells.Item(125, 2).Value = 68.599999999999994
$ws.Cells.Item(126, 2).Value = 61.4
$ws.Cells.Item(127, 2).Value = 64.3
$ws.Cells.Item(128, 2).Value = 67.099999999999994
$ws.Cells.Item(129, 2).Value = 61.4
$ws.Cells.Item(130, 2).Value = 58.6
$ws.Cells.Item(131, 2).Value = 64.3
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Values = $ws.Range("B5:B131")
Write-Output $ser.Formula
